$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cells value while preserving it as text, even when the
# string looks like a plain number (e.g. "24.60"), and without leaving a
# residual explicit cell style behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Update price (D) and volume (E) columns for rows with changed values
$ws.Range("D2").Value = '27.513.07'
$ws.Range("E2").Value = '  -2.94%  '

$ws.Range("D3").Value = '1.660.23'
$ws.Range("E3").Value = '  -3.80%  '

$ws.Range("E4").Value = '  -0.12%  '

Set-TextValue $ws.Range("D5") '214.38'
$ws.Range("E5").Value = '  -2.01%  '

$ws.Range("E6").Value = '  -2.54%  '

$ws.Range("E7").Value = '  -0.05%  '

Set-TextValue $ws.Range("D8") '24.60'
$ws.Range("E8").Value = '  +2.83%  '

Set-TextValue $ws.Range("D9") '0.264'
$ws.Range("E9").Value = '  -1.51%  '

Set-TextValue $ws.Range("D10") '0.0619'
$ws.Range("E10").Value = '  -2.55%  '

Set-TextValue $ws.Range("D11") '0.0878'
$ws.Range("E11").Value = '  -1.78%  '

$ws.Range("D12").Value = '1.893.49'
$ws.Range("E12").Value = '  -3.98%  '

$ws.Range("D13").Value = '1.657.88'
$ws.Range("E13").Value = '  -4.21%  '

$ws.Range("E14").Value = '  -2.43%  '

Set-TextValue $ws.Range("D15") '0.566'
$ws.Range("E15").Value = '  +0.42%  '

Set-TextValue $ws.Range("D16") '65.90'
$ws.Range("E16").Value = '  -2.46%  '

$ws.Range("D17").Value = '27.520.44'
$ws.Range("E17").Value = '  -2.82%  '

Set-TextValue $ws.Range("D18") '240.96'
$ws.Range("E18").Value = '  -1.92%  '

$ws.Range("E19").Value = '  -2.89%  '

Set-TextValue $ws.Range("D20") '7.61'
$ws.Range("E20").Value = '  -3.64%  '

$ws.Range("E21").Value = '  +0.07%  '

Set-TextValue $ws.Range("D22") '4.45'
$ws.Range("E22").Value = '  -3.49%  '

Set-TextValue $ws.Range("D23") '9.41'
$ws.Range("E23").Value = '  -2.36%  '

$ws.Range("E24").Value = '  -1.30%  '

Set-TextValue $ws.Range("D25") '145.99'
$ws.Range("E25").Value = '  -2.22%  '

Set-TextValue $ws.Range("D26") '7.22'
$ws.Range("E26").Value = '  -2.78%  '

Set-TextValue $ws.Range("D27") '16.29'
$ws.Range("E27").Value = '  -1.80%  '

$ws.Range("E28").Value = '  -0.17%  '

$ws.Range("E29").Value = '  -2.14%  '

Set-TextValue $ws.Range("D30") '0.0500'
$ws.Range("E30").Value = '  -3.14%  '

$ws.Range("E31").Value = '  -0.93%  '

$ws.Range("E32").Value = '  -2.92%  '

$ws.Range("D33").Value = '1.457.29'
$ws.Range("E33").Value = '  -1.70%  '

$ws.Range("E34").Value = '  -4.61%  '

$ws.Range("E35").Value = '  -4.23%  '

$ws.Range("E36").Value = '  -1.13%  '

Set-TextValue $ws.Range("D37") '0.925'
$ws.Range("E37").Value = '  -5.08%  '

$ws.Range("E38").Value = '  -4.60%  '

$ws.Range("E39").Value = '  -2.82%  '

$ws.Range("E40").Value = '  -0.56%  '

Set-TextValue $ws.Range("D41") '0.999'
$ws.Range("E41").Value = '  -0.08%  '

Set-TextValue $ws.Range("D42") '66.78'
$ws.Range("E42").Value = '  -4.09%  '

$ws.Range("E43").Value = '  -3.46%  '

$ws.Range("E44").Value = '  -2.77%  '

Set-TextValue $ws.Range("D47") '1.71'
$ws.Range("E47").Value = '  -0.45%  '

Set-TextValue $ws.Range("D48") '88.54'
$ws.Range("E48").Value = '  -1.84%  '

$ws.Range("E49").Value = '  -6.39%  '

$ws.Range("E50").Value = '  -1.59%  '

Set-TextValue $ws.Range("D51") '7.83'
$ws.Range("E51").Value = '  -3.65%  '

# Row 45 and 46: coin swap (RocketPoolETH and TrustWalletToken swap rows) with updated price/volume
$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D45").Value = '1.802.44'
$ws.Range("E45").Value = '  -4.00%  '

$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D46") '0.788'
$ws.Range("E46").Value = '  -2.13%  '
